$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab 9")

# --- Row 10 (Fortigate Firewall 1101E 7.0.8) ---
$ws.Range("F10").Value = "https://support.fortinet.com/Information/ProductLifeCycle.aspx"
$ws.Range("L10").Value = "7.0.8"
$ws.Range("M10").Value = 44847
$ws.Range("AA10").ClearContents()

# --- Row 13 (Fortinet/FortiManager-2000E 7.0.2) ---
$ws.Range("F13").Value = "https://docs.fortinet.com/document/fortimanager/7.0.5/release-notes/441895/change-log"

# M13 needs to become literal text "2022-10-13" while keeping its existing
# date number-format/style (not auto-converted back into a date serial).
# Temporarily switch to a text format so the assignment is stored as a
# string, then restore the original date format code - this keeps the
# same underlying style (no new style index / no quote-prefix) since the
# cell's content is already text by the time the format is restored.
$ws.Range("M13").NumberFormat = "@"
$ws.Range("M13").Value = "2022-10-13"
$ws.Range("M13").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("AA13").Value = "2023-01-30 19:25:55 (bot) : Updated by bot at this time"
# AA is outside the defined table (A1:Z23) and has no explicit column
# style in the source file; writing into a previously-empty cell here
# otherwise inherits the sheet's column-27 default style, so reset it
# back to the workbook's default "Normal" style to match.
$ws.Range("AA13").Style = "Normal"
